$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the F:AU price series for rows 4, 11, 13, 20, 21, 23, 25
# (historic + future scenario values replacing the placeholder zeros,
#  and rows 20/21 swapping which one carries the series).
$rowValues = @(0.06371098956,0.083524673699999996,0.064409172048000007,0.057368552832,0.053267649372,0.041721915743999999,0.043926702588000001,0.06210149518799999,0.065041210944000005,0.064503506015999995,0.038893990464000001,0.06078174436799999,0.058999937831999999,0.057270847428000003,0.055594473227999999,0.053970815159999998,0.054139506912,0.054308198664000003,0.054476890416000005,0.054645582132000006,0.054814273883999995,0.054824817131999999,0.054845903591999999,0.054856446840000003,0.054866990051999999,0.054877533300000003,0.055151657388,0.055436324687999999,0.055720992024,0.055995116111999997,0.056279783412000003,0.056743685711999996,0.057197044800000006,0.057671490311999994,0.058135392612,0.058609838159999997,0.059084283707999993,0.059569272468000006,0.060054261228000005,0.060539249987999998,0.061034781996000008,0.061034781996000008)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(4, 6 + $i).Value = $rowValues[$i]
}

$rowValues = @(0.25468054685999997,0.25118802651600003,0.24769550620800002,0.24420298590000003,0.24071046555599998,0.237217945248,0.233725424904,0.23023290459600002,0.22674038425199999,0.22986783781199999,0.134219882628,0.22270075671600001,0.21918237145199998,0.21566398618800001,0.212145600888,0.20849690505599999,0.21605491785600001,0.22413417292799998,0.232083117432,0.24016237250399999,0.24850224871200002,0.25162970230800003,0.25475715586799996,0.25749367776000004,0.26062113131999998,0.263618274348,0.26583355396800001,0.26817914411999999,0.27039442374,0.27300063506399996,0.27521591468400003,0.27834336824400002,0.28094957956800004,0.28407703312799998,0.28720448672399995,0.290331940284,0.29345939387999997,0.29658684743999997,0.29971430103599994,0.30284175459599999,0.30636013985999999,0.30636013985999999)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(11, 6 + $i).Value = $rowValues[$i]
}

$rowValues = @(0.33002762695199994,0.33002762695199994,0.33002762695199994,0.33002762695199994,0.33002762695199994,0.33002762695199994,0.33002762695199994,0.33002762695199994,0.33002762695199994,0.334579732164,0.3276,0.32414782438799999,0.31902670605600003,0.31390558768799998,0.30878446935600001,0.30347367994800001,0.29861807785200001,0.29357552590071445,0.28853297394942851,0.28349042199814295,0.27844787004685706,0.2734053180955715,0.26836276614428556,0.263320214193,0.25827766224171445,0.25323511029042856,0.24819255833914297,0.24315000638785703,0.23810745443657147,0.23306490248528555,0.228022350534,0.22297979858271444,0.21793724663142852,0.21289469468014297,0.20785214272885705,0.20280959077757149,0.19776703882628557,0.19272448687499999,0.18768193492371443,0.18263938297242852,0.17759683102114299,0.17759683102114299)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(13, 6 + $i).Value = $rowValues[$i]
}

$rowValues = @(0.074162839057389032,0.074162839057389032,0.074162839057389032,0.074162839057389032,0.074162839057389032,0.074162839057389032,0.074162839057389032,0.074162839057389032,0.074162839057389032,0.074162839057389032,0.059270793267966984,0.086844213682472068,0.09987405364380024,0.089848186169571506,0.080709805803888562,0.074545217103161129,0.06968777153502935,0.065984808790841606,0.064308412595709019,0.0632016756145867,0.062094938633464367,0.061968961426494336,0.061842984219524298,0.061717007012554274,0.061591029805584237,0.061465052598614199,0.061339075391644175,0.061213098184674145,0.061087120977704121,0.060961143770734083,0.060835166563764059,0.060709189356794022,0.060583212149823984,0.06045723494285396,0.060331257735883922,0.060205280528913899,0.060079303321943882,0.059953326114973844,0.059827348908003807,0.059701371701033769,0.059575394494063731,0.059575394494063731)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(20, 6 + $i).Value = $rowValues[$i]
}

$rowValues = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(21, 6 + $i).Value = $rowValues[$i]
}

$rowValues = @(0.30862373122800002,0.30439098280799998,0.30015823442399997,0.29592548600399998,0.29169273762000003,0.28745998919999999,0.28322724081600004,0.27899449239599999,0.27476174401199999,0.278551561176,0.16264631973599999,0.26986656351599997,0.26560301922000001,0.261339474924,0.25707593059199996,0.25265447726400003,0.26181320205600001,0.27160356304800004,0.281236015008,0.29102637599999998,0.30113255509199999,0.30492237222000002,0.30871218938400002,0.31202827941600003,0.31581809657999999,0.31945000467600004,0.32213445850799999,0.32497682137199996,0.32766127520400001,0.33081945616800001,0.33350390996399998,0.33729372712799999,0.34045190809199999,0.344241725256,0.34803154242000001,0.35182135958399996,0.35561117674800002,0.359400993876,0.36319081104000001,0.36698062824,0.37124417268000004,0.37124417268000004)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(23, 6 + $i).Value = $rowValues[$i]
}

$rowValues = @(0.54773038152000009,0.54773038152000009,0.54773038152000009,0.54773038152000009,0.54773038152000009,0.54773038152000009,0.54773038152000009,0.54773038152000009,0.54773038152000009,0.55528528320000003,0.54,0.53797196664000002,0.52947270215999998,0.52097343768000004,0.51247417319999999,0.50366012076,0.49613028768,0.48802579947000002,0.47992131125999998,0.47181682305,0.46371233483999996,0.45560784663000004,0.44750335842,0.43939887021000001,0.43129438199999998,0.42318989378999999,0.41508540557999996,0.40698091736999997,0.39887642915999999,0.39077194094999995,0.38266745274000002,0.37456296452999999,0.36645847632,0.35835398811000002,0.35024949989999998,0.34214501169,0.33404052347999996,0.32593603526999998,0.31783154705999994,0.30972705884999996,0.30162257064000003,0.29351808242999999)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(25, 6 + $i).Value = $rowValues[$i]
}

# Add a new (currently blank) summary row 28 formatted as whole numbers (#,##0),
# extending the sheet dimension to A1:AU28.
$ws.Range("F28:AU28").NumberFormat = "#,##0"

# Reset the view: select F28:AU29 with F28 active (mirrors the saved selection state).
$ws.Range("F28:AU29").Select() | Out-Null
